$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first data row (row 2). This removes "A 61099-2024" and
# shifts all subsequent rows up by one, matching the diff where every
# row's content becomes the content of the row that used to follow it,
# and the last row (old row 37) disappears.
$ws.Rows(2).Delete()

# The "Förändrad" (changed) date column C was recalculated/bumped from
# 45665 to 45666 for every remaining data row.
$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45665) {
        $cell.Value = 45666
    }
}
